$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: test_sql/0
$ws.Range("C2").Value = "['select count(*) from products']"
$ws.Range("D2").Value = "[True]"

# Row 3: test_sql/1
$ws.Range("C3").Value = "['select count(*) from orders']"
$ws.Range("D3").Value = "[True]"

# Row 4: test_sql/2
$ws.Range("B4").Value = "What is the average , minimum , and maximum price of all spanish products? | products : product_id, country, price"
$ws.Range("C4").Value = "[""select avg(price), min(price), max(price) from products where counstry = 'SPAIN'""]"

# Row 5: test_sql/3
$ws.Range("C5").Value = "['select country, count(*) from products group by country']"
$ws.Range("D5").Value = "[True]"

# Row 6: test_sql/4
$ws.Range("C6").Value = "[""select store, count(*) from sales where store = 'STORE1'""]"
$ws.Range("D6").Value = "[True]"

# Row 7: test_sql/5
$ws.Range("C7").Value = "['select name from customers where age > (select avg(age) from customers)']"
$ws.Range("D7").Value = "[True]"

# Row 8: test_sql/6
$ws.Range("C8").Value = "['select count(*) from orders as t1 join order_product as t2 on t1.order_id = t2.order_id where t2.price >= 100']"
$ws.Range("D8").Value = "[True]"

# Row 9: test_sql/7
$ws.Range("C9").Value = "['select store, count(store) from sales group by store order by count(store) desc limit 1']"
$ws.Range("D9").Value = "[True]"

# Row 10: test_sql/8
$ws.Range("C10").Value = "[""select t1.sale_id from sales as t1 join order_product as t2 on t1.product_id = t2.product_id where t2.country = 'SPAIN' and t1.quantity > 1""]"

# Row 11: test_sql/9
$ws.Range("C11").Value = "['select name from customers as t1 join sales as t2 on t1.customer_id = t2.customer_id group by t1.customer_id order by sum(t2.sales_id) desc limit 1']"
